$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.717.31'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '1.601.39'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = "'211.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("D10").Value = "'19.65"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("D11").Value = "'0.0843"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '1.825.99'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '1.617.74'
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("D16").Value = "'65.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '26.690.71'
$ws.Range("D18").Value = '0.0₃0744'
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("D19").Value = "'210.58"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("E20").Value = '  +2.31%  '
$ws.Range("D21").Value = "'1.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").Value = "'143.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").Value = "'7.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("D29").Value = "'15.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("D30").Value = "'0.0513"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").Value = "'1.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").Value = '1.297.09'
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").Value = "'0.609"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("E37").Value = '  +1.08%  '
$ws.Range("E38").Value = '  +20.33%  '
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").Value = "'63.25"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").Value = '1.737.03'
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = "'91.01"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("E47").Value = '  -2.04%  '
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("E51").Value = '  +0.02%  '
